# Generate Report for Handoff
#
# b.md has been handed off again: a fresh handoff package was generated for
# both the zh-cn and de-de targets. Update the localization-status report so
# that:
#   - the Overview sheet shows b.md's status as "Ready for handoff"
#   - each language sheet (zh-cn / de-de) shows b.md's status as
#     "Ready for handoff", with its "Latest Handoff File" / "Latest Handoff
#     Datetime" columns pointing at the newly generated handoff package.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"

# --- Overview sheet: row 3 is b.md -----------------------------------------
$overview.Range("B3").Value = $readyForHandoff
$overview.Range("C3").Value = $readyForHandoff

# --- zh-cn sheet: row 3 is b.md ---------------------------------------------
$zhcn.Range("B3").Value = $readyForHandoff
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-08 06:16:40"

# --- de-de sheet: row 3 is b.md ---------------------------------------------
$dede.Range("B3").Value = $readyForHandoff
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-08 06:16:42"
